$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-10 from 45183 to 45184
$ws.Range("C2:C10").Value = 45184
